$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 99, shifting existing
# rows 99-116 down to 101-118. The new rows inherit formatting (incl. the
# date-formatted style on column D) from the row above, same as Excel's
# native Insert behavior.
$ws.Range("A99:A100").EntireRow.Insert()

# Populate new row 99 (weekly Berenjena entry, Arica y Parinacota origin)
$ws.Cells.Item(99,1).Value2  = 6
$ws.Cells.Item(99,2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(99,3).Value2  = "Metropolitana"
$ws.Cells.Item(99,4).Value2  = 44504
$ws.Cells.Item(99,5).Value2  = 13
$ws.Cells.Item(99,6).Value2  = 100112001
$ws.Cells.Item(99,7).Value2  = "Berenjena"
$ws.Cells.Item(99,8).Value2  = "Sin especificar"
$ws.Cells.Item(99,9).Value2  = "Primera"
$ws.Cells.Item(99,10).Value2 = 400
$ws.Cells.Item(99,11).Value2 = 7000
$ws.Cells.Item(99,12).Value2 = 8000
$ws.Cells.Item(99,13).Value2 = 7575
$ws.Cells.Item(99,14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(99,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(99,16).Value2 = 152
$ws.Cells.Item(99,17).Value2 = 50
$ws.Cells.Item(99,18).Value2 = "Hortaliza"

# Populate new row 100 (weekly Berenjena entry, Provincia de Huasco origin)
$ws.Cells.Item(100,1).Value2  = 6
$ws.Cells.Item(100,2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(100,3).Value2  = "Metropolitana"
$ws.Cells.Item(100,4).Value2  = 44504
$ws.Cells.Item(100,5).Value2  = 13
$ws.Cells.Item(100,6).Value2  = 100112001
$ws.Cells.Item(100,7).Value2  = "Berenjena"
$ws.Cells.Item(100,8).Value2  = "Sin especificar"
$ws.Cells.Item(100,9).Value2  = "Primera"
$ws.Cells.Item(100,10).Value2 = 200
$ws.Cells.Item(100,11).Value2 = 11000
$ws.Cells.Item(100,12).Value2 = 13000
$ws.Cells.Item(100,13).Value2 = 12200
$ws.Cells.Item(100,14).Value2 = "$/caja 70 unidades"
$ws.Cells.Item(100,15).Value2 = "Provincia de Huasco"
$ws.Cells.Item(100,16).Value2 = 174
$ws.Cells.Item(100,17).Value2 = 70
$ws.Cells.Item(100,18).Value2 = "Hortaliza"
